# The commit swaps the contents of ppt/theme/theme1.xml (the design theme,
# used by the slide master -- originally the "Integral" palette) and
# ppt/theme/theme2.xml (the theme used only by the notes master --
# originally the stock "Office Theme" palette): after the edit, theme1.xml
# carries the Office Theme color values and theme2.xml carries the
# Integral color values.
#
# The PowerPoint object model only exposes one editable theme color
# scheme -- the design theme bound to the slide master (ppt/theme/theme1.xml)
# -- via Slide.ThemeColorScheme / Master.Theme.ThemeColorScheme /
# NotesMaster.Theme.ThemeColorScheme (they all resolve to the same
# underlying theme). So we drive the color change for theme1.xml through
# that surface, assigning it the 12 "Office" theme colors in
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $tcs.Item($i + 1).RGB = ToComRgb($officeColors[$i])
}
